$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text assignments (already safe as text: contain letters, URLs, % signs, multi-dot numbers, or special unicode digits) ---
$ws.Range("D2").Value = '69.318.75'
$ws.Range("E2").Value = '  +2.43%  '
$ws.Range("D3").Value = '3.388.97'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("E6").Value = '  +2.68%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("E9").Value = '  +8.15%  '
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("E11").Value = '  +3.50%  '
$ws.Range("E12").Value = '  +4.29%  '
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("E14").Value = '  +3.10%  '
$ws.Range("D15").Value = '3.934.32'
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").Value = '69.383.64'
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("D17").Value = '3.391.79'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("E20").Value = '  +2.24%  '
$ws.Range("E21").Value = '  +1.51%  '
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("E23").Value = '  +2.19%  '
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("E27").Value = '  +2.67%  '
$ws.Range("E28").Value = '  +2.66%  '
$ws.Range("E29").Value = '  +2.92%  '
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E32").Value = '  +10.96%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E33").Value = '  -2.92%  '
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("E35").Value = '  +2.29%  '
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = '3.669.55'
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("E38").Value = '  +5.29%  '
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("D40").Value = '0.0₃0722'
$ws.Range("E40").Value = '  +7.90%  '
$ws.Range("E41").Value = '  +3.67%  '
$ws.Range("E42").Value = '  +3.28%  '
$ws.Range("E43").Value = '  +1.93%  '
$ws.Range("E44").Value = '  +4.21%  '
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("E48").Value = '  +5.48%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +2.05%  '
$ws.Range("E51").Value = '  +4.85%  '

# --- Cells whose new value is a "pure" decimal number: force Text format first so Excel stores the literal
#     string (matching the source inlineStr cells) instead of silently converting to a Number, then clear
#     the temporary format so the cell keeps its original (default) style. ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.06'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.51'
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.59'
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000284'
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '679.51'
$ws.Range("D13").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.74'
$ws.Range("D19").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.23'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.28'
$ws.Range("D24").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.69'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.92'
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.79'
$ws.Range("D29").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.64'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '555.28'
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.107'
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.55'
$ws.Range("D35").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.56'
$ws.Range("D39").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.28'
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("D42").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0425'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.29'
$ws.Range("D45").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.40'
$ws.Range("D48").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.48'
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.67'
$ws.Range("D51").ClearFormats()
